# Update the "Totals" worksheet for the newly-completed month (row 17,
# 43862 = Jan-20) with the final figures now that the month has closed out
# (previously only a partial month was recorded).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Totals")

# Raw input values for row 17
$ws.Range("B17").Value = 2749
$ws.Range("C17").Value = 83167
$ws.Range("E17").Value = 45873
$ws.Range("F17").Value = 131789

# Match the number formatting used by the rest of the table: column B
# (rows 10-17) picks up the thousands-separated "#,##0" format that the
# earlier rows already had, and C17 picks up the same "Comma" formatting
# used by the other cells in column C.
$ws.Range("B10:B17").NumberFormat = "#,##0"
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection on the running-total cell, as in the saved file.
[void]$ws.Range("C30").Select()
